# Update column B values for specific rows as per data imputation re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 5.492999999999999
    27 = 6.55
    32 = 6.614
    36 = 8.704000000000001
    38 = 5.431
    46 = 6.62
    54 = 5.05
    55 = 4.684
    56 = 4.737
    67 = 4.988999999999999
    69 = 5.002
    72 = 5.380999999999999
    83 = 5.551
    86 = 5.146000000000001
    91 = 6.109
    93 = 5.380999999999999
    99 = 5.217000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
